# Generate Report for Handoff
# Adds a new "462266fd-..." file entry (status "Ready for handoff") as a new
# row on each of the three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

# ---- long, repeated literal fragments -------------------------------------------------
$oPad  = "ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo"
$newId   = "462266fd-499a-4cc0-b7b1-ecff9f616d55" + $oPad + ".md"
$newIdE2e = "e2e\" + $newId

$oPad2 = "oooooooooooooooooooooooooooooooooooooooo"
$newZh = "462266fd-499a-4cc0-b7b1-ecff9f616d55" + $oPad2 + ".6739222501bb04c098372b1353e10a300fdda172.zh-cn.xlf"
$newDe = "462266fd-499a-4cc0-b7b1-ecff9f616d55" + $oPad2 + ".6739222501bb04c098372b1353e10a300fdda172.de-de.xlf"

$hoDate   = "2016-08-16 22:26:31"
$zhDate   = "2016-08-16 22:26:26"
$deDate   = $hoDate
$status   = "Ready for handoff"

# =========================================================================
# Sheet: Overview
# =========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newId
$wsOverview.Range("B3").Value = $newIdE2e
$wsOverview.Range("B3").Style = "Hyperlink"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $hoDate
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0490dbe5923558288bf1d9580d42d3887332f9ea/e2e/" + $newId, "", "", $newIdE2e) | Out-Null

# =========================================================================
# Sheet: zh-cn
# =========================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = $newId
$wsZh.Range("A3").Style = "Hyperlink"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $status
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("G3").Value = $newZh
$wsZh.Range("H3").Value = $zhDate
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("O3").Value = "'False"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0490dbe5923558288bf1d9580d42d3887332f9ea/e2e/" + $newId, "", "", $newId) | Out-Null

# =========================================================================
# Sheet: de-de
# =========================================================================
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = $newId
$wsDe.Range("A3").Style = "Hyperlink"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $status
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("G3").Value = $newDe
$wsDe.Range("H3").Value = $deDate
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("O3").Value = "'False"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0490dbe5923558288bf1d9580d42d3887332f9ea/e2e/" + $newId, "", "", $newId) | Out-Null

# =========================================================================
# Column width tweaks (Overview cols E/F, zh-cn/de-de col C)
# =========================================================================
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797
$wsZh.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDe.Columns.Item(3).ColumnWidth = 17.2159881591797
